$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 12.49520910190846
$ws.Range("C2").Value = 9.288973803215931
$ws.Range("E2").Value = 11.63215065497509
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 3.59701488501062
$ws.Range("I2").Value = 16.5409470004262
$ws.Range("M2").Value = 14.45124636476535
$ws.Range("O2").Value = 17.08095375148841
$ws.Range("B3").Value = 11.80570860404007
$ws.Range("C3").Value = 8.815037415443363
$ws.Range("E3").Value = 11.56933059045996
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 3.599013418264768
$ws.Range("I3").Value = 16.70558851522055
$ws.Range("M3").Value = 14.12112120627908
$ws.Range("O3").Value = 17.21003898883817
$ws.Range("B4").Value = 11.36066876589177
$ws.Range("C4").Value = 8.509244669694702
$ws.Range("E4").Value = 11.53603841280217
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 3.600304064271541
$ws.Range("I4").Value = 16.81194457957259
$ws.Range("M4").Value = 13.91701508939771
$ws.Range("O4").Value = 17.29578512175504
$ws.Range("B5").Value = 11.17399413870485
$ws.Range("C5").Value = 8.380998531011246
$ws.Range("E5").Value = 11.52380862720955
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 3.60084604220815
$ws.Range("I5").Value = 16.85661156403726
$ws.Range("M5").Value = 13.83360514104817
$ws.Range("O5").Value = 17.3323495257964
$ws.Range("B6").Value = 11.14268032395631
$ws.Range("C6").Value = 8.359486800633196
$ws.Range("E6").Value = 11.52185885464774
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 3.600937006883586
$ws.Range("I6").Value = 16.86410863077479
$ws.Range("M6").Value = 13.81974434432213
$ws.Range("O6").Value = 17.33851872758876
$ws.Range("B7").Value = 11.35817253657252
$ws.Range("C7").Value = 8.507529680537861
$ws.Range("E7").Value = 11.53586805365167
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 3.600311308604816
$ws.Range("I7").Value = 16.81254160248511
$ws.Range("M7").Value = 13.91589098510911
$ws.Range("O7").Value = 17.29627168602941
$ws.Range("B8").Value = 12.26204256480773
$ws.Range("C8").Value = 9.128678145283994
$ws.Range("E8").Value = 11.60940106454469
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 3.597690823197618
$ws.Range("I8").Value = 16.59662371362445
$ws.Range("M8").Value = 14.33778146717664
$ws.Range("O8").Value = 17.12411038872115
$ws.Range("B9").Value = 13.85782499500589
$ws.Range("C9").Value = 10.22632494495438
$ws.Range("E9").Value = 11.79493573430996
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.593053831300742
$ws.Range("I9").Value = 16.21489919342875
$ws.Range("M9").Value = 15.14894406369789
$ws.Range("O9").Value = 16.8383891345102
$ws.Range("B10").Value = 14.91775964164868
$ws.Range("C10").Value = 10.95618605138706
$ws.Range("E10").Value = 11.95551373899527
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.589949588039928
$ws.Range("I10").Value = 15.95973689652087
$ws.Range("M10").Value = 15.72875019998894
$ws.Range("O10").Value = 16.66067405784695
$ws.Range("B11").Value = 15.3748541700418
$ws.Range("C11").Value = 11.27113094239676
$ws.Range("E11").Value = 12.0335817979061
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.588602372444955
$ws.Range("I11").Value = 15.84912405686558
$ws.Range("M11").Value = 15.9877997420632
$ws.Range("O11").Value = 16.5869455876722
$ws.Range("B12").Value = 15.54429826985264
$ws.Range("C12").Value = 11.38790935590368
$ws.Range("E12").Value = 12.06384158219008
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.588101498476572
$ws.Range("I12").Value = 15.80802175689236
$ws.Range("M12").Value = 16.08512611902366
$ws.Range("O12").Value = 16.56006034909027
$ws.Range("B13").Value = 15.50796824349323
$ws.Range("C13").Value = 11.36286992650414
$ws.Range("E13").Value = 12.05729400796989
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.58820895835363
$ws.Range("I13").Value = 15.81683900166775
$ws.Range("M13").Value = 16.06420080209343
$ws.Range("O13").Value = 16.56580439639374
$ws.Range("B14").Value = 15.38886774726727
$ws.Range("C14").Value = 11.28078830823281
$ws.Range("E14").Value = 12.03605746658612
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.588560979384723
$ws.Range("I14").Value = 15.84572683104499
$ws.Range("M14").Value = 15.99582278430733
$ws.Range("O14").Value = 16.58471293388638
$ws.Range("B15").Value = 15.31543914832705
$ws.Range("C15").Value = 11.23018664891812
$ws.Range("E15").Value = 12.02313948564294
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.588777810494804
$ws.Range("I15").Value = 15.86352358259543
$ws.Range("M15").Value = 15.95383634032871
$ws.Range("O15").Value = 16.59642996238399
$ws.Range("B16").Value = 14.88737904487002
$ws.Range("C16").Value = 10.93525745317024
$ws.Range("E16").Value = 11.95051075396168
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.590038933988595
$ws.Range("I16").Value = 15.96707552066158
$ws.Range("M16").Value = 15.71171786210776
$ws.Range("O16").Value = 16.66563649666398
$ws.Range("B17").Value = 14.6183225294293
$ws.Range("C17").Value = 10.74993216341455
$ws.Range("E17").Value = 11.90722419205379
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 3.590829185697984
$ws.Range("I17").Value = 16.0319994058729
$ws.Range("M17").Value = 15.56191193967058
$ws.Range("O17").Value = 16.70992257712131
$ws.Range("B18").Value = 14.46120992821949
$ws.Range("C18").Value = 10.64173228721237
$ws.Range("E18").Value = 11.88280106305523
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 3.591289831035013
$ws.Range("I18").Value = 16.0698560470923
$ws.Range("M18").Value = 15.47530961436611
$ws.Range("O18").Value = 16.73606379118818
$ws.Range("B19").Value = 14.40761069066706
$ws.Range("C19").Value = 10.60482289931864
$ws.Range("E19").Value = 11.87461397116346
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 3.591446849094076
$ws.Range("I19").Value = 16.08276199250974
$ws.Range("M19").Value = 15.44591529245984
$ws.Range("O19").Value = 16.74502933083618
$ws.Range("B20").Value = 14.64720847748285
$ws.Range("C20").Value = 10.76982677177537
$ws.Range("E20").Value = 11.91178322047426
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 3.590744429684367
$ws.Range("I20").Value = 16.02503494368055
$ws.Range("M20").Value = 15.57790505561378
$ws.Range("O20").Value = 16.70513892661902
$ws.Range("B21").Value = 15.42394973102658
$ws.Range("C21").Value = 11.30496529594025
$ws.Range("E21").Value = 12.04227644173069
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.588457330600867
$ws.Range("I21").Value = 15.83722049479701
$ws.Range("M21").Value = 16.01592867363513
$ws.Range("O21").Value = 16.57913088369899
$ws.Range("B22").Value = 15.91032089682995
$ws.Range("C22").Value = 11.64021919312693
$ws.Range("E22").Value = 12.13161173658671
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.587016690756077
$ws.Range("I22").Value = 15.71904450256035
$ws.Range("M22").Value = 16.29767551051517
$ws.Range("O22").Value = 16.50281025157077
$ws.Range("B23").Value = 15.65269335704186
$ws.Range("C23").Value = 11.46262176490143
$ws.Range("E23").Value = 12.08356988118979
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 3.587780651950418
$ws.Range("I23").Value = 15.78169924991075
$ws.Range("M23").Value = 16.14774515357761
$ws.Range("O23").Value = 16.54298833116532
$ws.Range("B24").Value = 14.6341567011177
$ws.Range("C24").Value = 10.76083756864011
$ws.Range("E24").Value = 11.9097206400011
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.590782728190531
$ws.Range("I24").Value = 16.0281819227537
$ws.Range("M24").Value = 15.57067604619131
$ws.Range("O24").Value = 16.70729949589663
$ws.Range("B25").Value = 13.44560244050117
$ws.Range("C25").Value = 9.942644885725178
$ws.Range("E25").Value = 11.74040447154775
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 3.594254888672898
$ws.Range("I25").Value = 16.31371591405754
$ws.Range("M25").Value = 14.93190724764554
$ws.Range("O25").Value = 16.91007090112527